# HBIS Linelist workbook update
# - Rename Formname "frmnumerictwo" -> "FrmDataID2" (row 8, Sheet1)
# - Rename Tablename "tblMainques" -> "tblLinelist " (rows 2-10, Sheet1, column D)
# - Rename SQL target tables in generated INSERT statements:
#     tblQuestion  -> tblQuestionLList   (Sheet1, column U)
#     tblOptions   -> tblOptionsLList    (Sheet2, column H)
# - Drop the stray vertical-top alignment on Sheet1!C8 (keeps wrap only)
# - Move the active selection on Sheet1 to E16

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: Formname for row 8 changes first so the new shared string
#     ("FrmDataID2") is interned before "tblLinelist " -- matches author order.
$ws1.Cells.Item(8, 3).Value = "FrmDataID2"

# --- Sheet1: Tablename column (D) for data rows 2-10
for ($r = 2; $r -le 10; $r++) {
    $ws1.Cells.Item($r, 4).Value = "tblLinelist "
}

# --- Sheet1: rewrite the generated INSERT formulas (column U, rows 2-10)
for ($r = 2; $r -le 10; $r++) {
    $cell = $ws1.Cells.Item($r, 21)
    $f = $cell.Formula
    $cell.Formula = $f.Replace("tblQuestion (", "tblQuestionLList (")
}

# --- Sheet2: rewrite the generated INSERT formulas (column H, rows 2-5)
for ($r = 2; $r -le 5; $r++) {
    $cell = $ws2.Cells.Item($r, 8)
    $f = $cell.Formula
    $cell.Formula = $f.Replace("tblOptions (", "tblOptionsLList (")
}

# --- Sheet1: C8 loses its "vertical: top" alignment, keeping wrap text only
$ws1.Cells.Item(8, 3).VerticalAlignment = -4107

# --- Sheet1: move the live selection to E16
$ws1.Range("E16").Select()

Write-Output "edit complete"
